# Writing up background of CD theory
# ------------------------------------------------------------
# 1. Rename Sheet1 -> Programming
# 2. Duplicate it to create "Writing" and "Structure" sheets
#    (duplicating preserves conditional formatting / dxf / page
#    setup so the new sheets inherit the same look as Programming)
# 3. Trim each duplicate down to the rows it actually needs and
#    fill in the new task/structure text
# 4. Append two new tasks to the bottom of Programming
# 5. Leave "Writing" as the active sheet/tab
# ------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsProgramming = $wb.Worksheets.Item(1)
$wsProgramming.Name = "Programming"

# ---- Programming: two new rows at the bottom -------------------------
$wsProgramming.Range("B22").Value = 1

$wsProgramming.Range("A9").Copy()
$wsProgramming.Range("A23").PasteSpecial(-4122)
$wsProgramming.Range("A23").Value = "Handle radius changes"
$wsProgramming.Range("B23").Value = 1

$wsProgramming.Range("A10").Copy()
$wsProgramming.Range("A24").PasteSpecial(-4122)
$wsProgramming.Range("A24").Value = "Add attribute outcome to CD events"

$wsProgramming.Range("A13").Select() | Out-Null

# ---- Writing sheet -----------------------------------------------------
$wsProgramming.Copy($null, $wsProgramming)
$wsWriting = $wb.Worksheets.Item(2)
$wsWriting.Name = "Writing"

$wsWriting.Range("A5:B24").EntireRow.Delete()
$wsWriting.Range("A2:B4").ClearContents()
$wsWriting.Cells.Validation.Delete()

$fcWriting = $wsWriting.Range("A2:A100").FormatConditions.Item(1)
$fcWriting.ModifyAppliesToRange($wsWriting.Range("A2:A1000"))

$wsWriting.Columns("A").ColumnWidth = 29.1
$wsWriting.Range("A2:A4").Style = "Normal"

$wsWriting.Range("A2").Value = "Read example dissertations"
$wsWriting.Range("B2").Value = 1
$wsWriting.Range("A3").Value = "Read papers on CD/Metaphor stuff"
$wsWriting.Range("A4").Value = "Write more detail on implementation"

$wsWriting.Range("B4").Select() | Out-Null

# ---- Structure sheet ----------------------------------------------------
$wsProgramming.Copy($null, $wsWriting)
$wsStructure = $wb.Worksheets.Item(3)
$wsStructure.Name = "Structure"

$wsStructure.Range("A7:B24").EntireRow.Delete()
$wsStructure.Cells.Validation.Delete()
$wsStructure.Range("A2:A100").FormatConditions.Delete()
$wsStructure.Range("A1:B6").Clear()

$wsStructure.Columns("A").ColumnWidth = 22.42

$wsStructure.Range("A1").Value = "Introduction"
$wsStructure.Range("A2").Value = "Background"
$wsStructure.Range("A3").Value = "Literature Review"
$wsStructure.Range("A4").Value = "Implementation"
$wsStructure.Range("A5").Value = "Results"
$wsStructure.Range("A6").Value = "Conclusion + Future Work"

$wsStructure.Range("A7").Select() | Out-Null

# ---- Leave "Writing" as the active tab ----------------------------------
$wsWriting.Activate()
$wsWriting.Range("B4").Select() | Out-Null
